$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A26").Value = "isPartOfBookChapter"
$ws.Range("D26").Value = "Appartient au chapitre"
$ws.Range("C26").Value = "Teil des Kapitels"
$ws.Range("B26").Value = "Part of chapter"
$ws.Range("E26").Value = "Appartiene al capitolo"
$ws.Range("G26").Value = "Belongs to following chapter"
$ws.Range("H26").Value = "Gehört zu folgendem Kapitel"
$ws.Range("I26").Value = "Appartient au chapitre suivant"
$ws.Range("J26").Value = "Appartiene al seguente capitolo"
$ws.Range("M26").Value = ":BookChapter"
